# DOMA-2542 Localization for Excel template (contacts)
# Replace hard-coded Russian sheet name and header labels with i18n placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Localize the sheet (tab) name.
$ws.Name = "{d.i18n.sheetName}"

# Localize the header row labels.
$ws.Range("A1").Value = "{d.i18n.name}"
$ws.Range("B1").Value = "{d.i18n.address}"
$ws.Range("C1").Value = "{d.i18n.unitName}"
$ws.Range("D1").Value = "{d.i18n.phone}"
$ws.Range("E1").Value = "{d.i18n.email}"
